$d = $word.ActiveDocument

function Merge-ParagraphRuns($index) {
    # Collapse a paragraph's runs (which may be split up by proofErr-marked
    # spell/grammar-check fragments) back into a single run by replacing the
    # paragraph's own text with itself via Find/Replace (scoped to just that
    # paragraph's Range so we never touch other, identical-looking text
    # elsewhere in the document).
    $p = $d.Paragraphs.Item($index)
    $full = $p.Range.Text
    $full = $full.TrimEnd([char]13, [char]7)
    $rng = $p.Range
    $find = $rng.Find
    $find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, $full, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Tidy up the proofing-split runs (an/SQLite/etc.) throughout the intro
#    and task-list paragraphs -- this merges each paragraph's text back into
#    a single run with no leftover <w:proofErr/> markers, matching the
#    "cleaned up" wording from the diff.
# ---------------------------------------------------------------------------
Merge-ParagraphRuns 2
Merge-ParagraphRuns 3
Merge-ParagraphRuns 4
Merge-ParagraphRuns 8
Merge-ParagraphRuns 9
Merge-ParagraphRuns 11

# ---------------------------------------------------------------------------
# 2) Add the three new sub-bullets (David / Yusuf / Anthony) right before the
#    "Course information", "Add tasks", and "Grade calculation" user-story
#    bullets respectively, at the nested (ilvl=1) list level.
# ---------------------------------------------------------------------------
function Insert-SubBullet($beforeText, $name) {
    $d2 = $word.ActiveDocument
    for ($i = 1; $i -le $d2.Paragraphs.Count; $i++) {
        $p = $d2.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($beforeText)) {
            $p.Range.InsertParagraphBefore()
            $newp = $d2.Paragraphs.Item($i)
            $newp.Range.ListFormat.ListLevelNumber = 2
            $newp.Range.InsertBefore($name)
            break
        }
    }
}

Insert-SubBullet "Course information" "David"
Insert-SubBullet "Add tasks" "Yusuf"
Insert-SubBullet "Grade calculation" "Anthony"

# ---------------------------------------------------------------------------
# 3) Final paragraph: merge its proofing-split runs, then attribute the
#    SQLite-database work to "(Ian & Graham)" by splitting the run right
#    after "We ".
# ---------------------------------------------------------------------------
$lastIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("We are also implementing")) {
        $lastIdx = $i
        break
    }
}

Merge-ParagraphRuns $lastIdx

$p = $d.Paragraphs.Item($lastIdx)
$full = $p.Range
$insPos = $full.Start + 3
$insertRange = $d.Range($insPos, $insPos)
$insertRange.InsertAfter("(Ian & Graham) ")

# Force the newly-typed text into its own run (rather than silently
# re-merging with its neighbours) by toggling a character property on it.
$splitRange = $d.Range($insPos, $insPos + 15)
$splitRange.Bold = 1
$splitRange.Bold = 0

Write-Host "Done"
